$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet methodID -> methodID2
$ws.Name = "methodID2"

# Insert two new columns before column G (keeps existing G..K data/styles,
# shifting them to I..M) to make room for the new "thickness" / "scanningTime"
# columns.
$ws.Columns("G:H").Insert()

# --- Group header row (row 1) ---
$ws.Range("B1").Value = "instructionID2"
$ws.Range("G1").Value = "sample002"
$ws.Range("H1").Value = "condition002"

# --- Column header row (row 2) ---
$ws.Range("B2").Value = "DATE"
$ws.Range("G2").Value = "thickness"
$ws.Range("H2").Value = "scanningTime"
$ws.Range("I2").Value = "comments"

# --- Data rows: new "thickness" (text, quote-prefixed) and "scanningTime"
#     (numeric) columns ---
$ws.Range("G3").Value = "'0.10"
$ws.Range("H3").Value = 10

$ws.Range("G4").Value = "'0.12"
$ws.Range("H4").Value = 12

$ws.Range("G5").Value = "'0.11"

$ws.Range("G6").Value = "'0.10"
$ws.Range("H6").Value = 14

$ws.Range("G7").Value = "'0.12"
$ws.Range("H7").Value = 10

$ws.Range("G8").Value = "'0.10"

$ws.Range("G9").Value = "'0.10"
$ws.Range("H9").Value = 20
